{"js": "// Fixed case-insensitive keyword matching\n// The feature bullet \"Supports custom keyword mappings through a configuration\n// file.\" should read \"Supports case-insensitive custom keyword mappings\n// through a configuration file.\"\n\nconst results = context.document.body.search(\"Supports custom keyword mappings through a configuration file.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target sentence to update.\");\n}\n\nresults.items[0].insertText(\n  \"Supports case-insensitive custom keyword mappings through a configuration file.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Fixed case-insensitive keyword matching\n# The feature bullet \"Supports custom keyword mappings through a configuration\n# file.\" should read \"Supports case-insensitive custom keyword mappings\n# through a configuration file.\"\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n\n$found = $rng.Find.Execute(\n    \"Supports custom keyword mappings through a configuration file.\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"Supports case-insensitive custom keyword mappings through a configuration file.\",\n    2\n)\n\nif (-not $found) {\n    throw \"Could not find target sentence to update.\"\n}\n"}
